# Adds a new column D ("counting-in-the-amazon") with Good/No values,
# mirroring columns B and C already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell D1 - copy formatting from C1 (bold header style) then set value.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D1").Value = "counting-in-the-amazon"

$values = @(
    "Good", # row 2
    "Good", # row 3
    "No",   # row 4
    "Good", # row 5
    "Good", # row 6
    "No",   # row 7
    "No",   # row 8
    "Good", # row 9
    "Good", # row 10
    "No",   # row 11
    "Good", # row 12
    "Good", # row 13
    "No",   # row 14
    "No",   # row 15
    "No",   # row 16
    "Good", # row 17
    "Good", # row 18
    "Good", # row 19
    "Good", # row 20
    "Good", # row 21
    "No",   # row 22
    "No",   # row 23
    "Good", # row 24
    "No",   # row 25
    "No"    # row 26
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 4).Value = $v
    $row = $row + 1
}
